{"js": "// Update the nomination-form merge placeholders:\n//  1. Append the \"old company name\" conditional block after {companyName}.\n//  2. Rename every {shareholderName_N} placeholder to\n//     {shareholderNameCertificate_N} (the combined intro line plus the\n//     three individual \"Name of the Nth Security Holder(s)\" lines).\n\nconst body = context.document.body;\n\nasync function replaceExact(searchText, newText) {\n  const results = body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n  return results.items.length;\n}\n\n// 1) {companyName} -> {companyName} {#hasCompanyOldName}[{companyOldName}]{/hasCompanyOldName}\nawait replaceExact(\n  \"{companyName}\",\n  \"{companyName} {#hasCompanyOldName}[{companyOldName}]{/hasCompanyOldName}\"\n);\n\n// 2) Combined intro line listing up to three shareholders.\nawait replaceExact(\n  \"{#hasShareholder_1}{shareholderName_1}{/hasShareholder_1}\" +\n    \"{#hasShareholder_2}; {shareholderName_2}{/hasShareholder_2}\" +\n    \"{#hasShareholder_3}; {shareholderName_3}{/hasShareholder_3}\",\n  \"{#hasShareholder_1}{shareholderNameCertificate_1}{/hasShareholder_1}\" +\n    \"{#hasShareholder_2}; {shareholderNameCertificate_2}{/hasShareholder_2}\" +\n    \"{#hasShareholder_3}; {shareholderNameCertificate_3}{/hasShareholder_3}\"\n);\n\n// 3) The three standalone \"Name of the Nth Security Holder(s)\" placeholders.\nfor (const n of [1, 2, 3]) {\n  await replaceExact(\n    `{#hasShareholder_${n}}{shareholderName_${n}}{/hasShareholder_${n}}`,\n    `{#hasShareholder_${n}}{shareholderNameCertificate_${n}}{/hasShareholder_${n}}`\n  );\n}\n", "ps1": "# Update the nomination-form merge placeholders:\n#  1. Append the \"old company name\" conditional block after {companyName}.\n#  2. Rename every {shareholderName_N} placeholder to\n#     {shareholderNameCertificate_N} (the combined intro line plus the\n#     three individual \"Name of the Nth Security Holder(s)\" lines).\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) {companyName} -> {companyName} {#hasCompanyOldName}[{companyOldName}]{/hasCompanyOldName}\nReplace-ExactText \"{companyName}\" \"{companyName} {#hasCompanyOldName}[{companyOldName}]{/hasCompanyOldName}\"\n\n# 2) Combined intro line listing up to three shareholders.\n$oldCombined = \"{#hasShareholder_1}{shareholderName_1}{/hasShareholder_1}{#hasShareholder_2}; {shareholderName_2}{/hasShareholder_2}{#hasShareholder_3}; {shareholderName_3}{/hasShareholder_3}\"\n$newCombined = \"{#hasShareholder_1}{shareholderNameCertificate_1}{/hasShareholder_1}{#hasShareholder_2}; {shareholderNameCertificate_2}{/hasShareholder_2}{#hasShareholder_3}; {shareholderNameCertificate_3}{/hasShareholder_3}\"\nReplace-ExactText $oldCombined $newCombined\n\n# 3) The three standalone \"Name of the Nth Security Holder(s)\" placeholders.\nfor ($n = 1; $n -le 3; $n++) {\n    $oldText = \"{#hasShareholder_$n}{shareholderName_$n}{/hasShareholder_$n}\"\n    $newText = \"{#hasShareholder_$n}{shareholderNameCertificate_$n}{/hasShareholder_$n}\"\n    Replace-ExactText $oldText $newText\n}\n"}
